$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: fill in F4, J4, N4 with values, matching the red style already used by D4/B4/C4
$ws.Range("F4").Font.Color = $ws.Range("D4").Font.Color
$ws.Range("F4").Value = 2862

$ws.Range("J4").Font.Color = $ws.Range("D4").Font.Color
$ws.Range("J4").Value = 3199

$ws.Range("N4").Font.Color = $ws.Range("D4").Font.Color
$ws.Range("N4").Value = 1361

# Row 6: add F6, J6, N6 with values, matching the red style used by D6
$ws.Range("F6").Font.Color = $ws.Range("D6").Font.Color
$ws.Range("F6").Value = 40136320

$ws.Range("J6").Font.Color = $ws.Range("D6").Font.Color
$ws.Range("J6").Value = 44043805

$ws.Range("N6").Font.Color = $ws.Range("D6").Font.Color
$ws.Range("N6").Value = 14458753

# Row 16: H16 gets a value with the same red style as G16
$ws.Range("H16").Font.Color = $ws.Range("G16").Font.Color
$ws.Range("H16").Value = 52884

# Row 17: H17 already has the blue family-3 font; just recolor to red and set the number
$ws.Range("H17").Font.Color = $ws.Range("D4").Font.Color
$ws.Range("H17").Value = 14105

# Row 18: add H18 with a value matching the red style used by G18
$ws.Range("H18").Font.Color = $ws.Range("G18").Font.Color
$ws.Range("H18").Value = 36861622

# The newly populated columns (F, H, J, N) get an auto-fit-style custom width,
# matching the existing best-fit width already used by column D.
$ws.Columns("F").ColumnWidth = $ws.Columns("D").ColumnWidth
$ws.Columns("H").ColumnWidth = $ws.Columns("D").ColumnWidth
$ws.Columns("J").ColumnWidth = $ws.Columns("D").ColumnWidth
$ws.Columns("N").ColumnWidth = $ws.Columns("D").ColumnWidth

# Update the selected cell shown when the sheet is opened
$ws.Range("N4").Select()
